$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("241002")

$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("D11").Value = 1
